$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arbeitsmatrix")

# Insert a new row at 136 (existing row 136 and below shift down by one)
$ws.Rows.Item(136).Insert()

# Fill in the new row 136 with the new time-entry data
$ws.Cells.Item(136, 1).Value = 22
$ws.Cells.Item(136, 2).Value = "Interface Design"
$ws.Cells.Item(136, 3).Value = "MockUps"
$ws.Cells.Item(136, 4).Value = "[FEATURE]"
$ws.Cells.Item(136, 5).Value = "Übersetzen, FreundRezeptAnsicht und Imports"
$ws.Cells.Item(136, 6).Value = 44464
$ws.Cells.Item(136, 7).Value = 44481
$ws.Range("I136").Formula = "=ROUNDUP(((SUM(K136-J136)*24*60/60)/0.25),0)*0.25"
$ws.Cells.Item(136, 10).Value = 0.55208333333333337
$ws.Cells.Item(136, 11).Value = 0.59375

# Insert another new row at 137 (blank separator row, but carries J/K formulas)
$ws.Rows.Item(137).Insert()
$ws.Range("A137").Clear()
$ws.Range("B137").Clear()
$ws.Range("C137").Clear()
$ws.Range("E137").Clear()
$ws.Range("I137").Clear()
$ws.Range("J137").Formula = "=K136"
$ws.Cells.Item(137, 11).Value = 0.70138888888888884

# Update selection/view to match target
$ws.Application.ActiveWindow.ScrollRow = 120
$ws.Range("L137").Select()
